$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RVL")

# "Global / DoSleep / millis / number / 2000" parameter rows (3 and 5) had their
# Param Value (column G) changed from "2000" to "500". Use a leading apostrophe
# so the numeric-looking value is stored as text (shared string), matching the
# original cell's string type (t="s").
$ws.Range("G3").Value = "'500"
$ws.Range("G5").Value = "'500"
